$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new entry "Parts for LED boards" in the next empty row (A8)
$ws.Range("A8").Value = "Parts for LED boards"

# Move selection to the next empty cell (A9), matching post-entry cursor position
$ws.Range("A9").Select()
